$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. New row 3: " " in A3 (centered like header), blank styled
#    cells for B3:H3 (same centered style), and I3 blank with
#    left/center alignment.
# ---------------------------------------------------------------
$ws.Range("A3").Value = " "

# Give A3:H3 the same look as the header row (font 標楷體 12,
# centered both ways) by copying the header's cell format.
$ws.Range("A2").Copy()
$ws.Range("A3:H3").PasteSpecial(-4122)

# I3: blank cell, left/center aligned, same font.
$ws.Range("I3").Font.Name = "標楷體"
$ws.Range("I3").Font.Size = 12
$ws.Range("I3").HorizontalAlignment = -4131
$ws.Range("I3").VerticalAlignment = -4108

# ---------------------------------------------------------------
# 2. New row 4: " " across A4:I4 with the plain body font
#    (標楷體 12, vertical-center only, no special border).
# ---------------------------------------------------------------
$ws.Range("A4:I4").Value = " "
$ws.Range("A4:I4").Font.Name = "標楷體"
$ws.Range("A4:I4").Font.Size = 12

# ---------------------------------------------------------------
# 3. Drop the row-2 "custom format" row default (ClearFormats
#    removes the row-level override) then restore the header's
#    visual formatting via a format-only paste from A1.
# ---------------------------------------------------------------
$ws.Rows(2).ClearFormats()
$ws.Range("A1").Copy()
$ws.Range("A2:I2").PasteSpecial(-4122)

# ---------------------------------------------------------------
# 4. Apply a thin box border (outline + inside) around the whole
#    used block A1:I3 to match the new bordered look.
# ---------------------------------------------------------------
$ws.Range("A1:I3").Borders.LineStyle = 1

# ---------------------------------------------------------------
# 5. Selection / dimension bookkeeping to match the saved file.
# ---------------------------------------------------------------
$ws.Range("I9").Select()
